$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86

# Columns whose values look numeric/date-like must be forced to Text so they
# are stored the same way as the rest of the sheet (inlineStr / shared-string
# "7756", "10/31/2025", "13", "810492559" rather than being auto-converted to
# numbers or date serials). Force text number format, assign, then clear the
# number-format stamp so no extra style index is introduced.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "A$row" "7756"
Set-TextValue "B$row" "10/31/2025"
$ws.Range("C$row").Value = "JURAMENTO 1415"
Set-TextValue "D$row" "13"
Set-TextValue "E$row" "810492559"
$ws.Range("F$row").Value = "NEW"
$ws.Range("G$row").Value = "Pendiente"
$ws.Range("H$row").Value = "Picada"
$ws.Range("I$row").Value = 1
$ws.Range("J$row").Value = "Cambio"
$ws.Range("K$row").Value = "Sin equipos"
$ws.Range("L$row").Value = "Pasante"
$ws.Range("M$row").Value = -58.446813
$ws.Range("N$row").Value = -34.556095
$ws.Range("O$row").Value = "Saavedra"
$ws.Range("P$row").Value = "Capital Norte"
$ws.Range("Q$row").Value = "BLO-A"
$ws.Range("R$row").Value = "Fuera de Poligono OVL"
